$d = $word.ActiveDocument

# Find the paragraph that holds "LOQ4037: Química Orgânica I (Requisito fraco)"
# and the one that holds the copyright notice; the text in between (an empty
# paragraph, the "Ver no Jupiter..." paragraph, and the copyright paragraph
# itself) must be removed, leaving the "LOQ4037..." paragraph directly
# followed by the blank paragraph that used to precede the final page break.

$count = $d.Paragraphs.Count
$idxLOQ = -1
$idxCopyright = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("LOQ4037: Química Orgânica I (Requisito fraco)")) { $idxLOQ = $i }
    if ($t.Contains("Powered by Jekyll and Github pages")) { $idxCopyright = $i }
}

if ($idxLOQ -gt 0 -and $idxCopyright -gt $idxLOQ) {
    $startRange = $d.Paragraphs.Item($idxLOQ + 1).Range
    $endRange = $d.Paragraphs.Item($idxCopyright).Range
    $delRange = $d.Range($startRange.Start, $endRange.End)
    $delRange.Delete()
}
